$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric CE test-case identifiers in column A (rows 25-38)
# with their "CE###" text labels, e.g. 100 -> "CE100".
$values = @(100, 110, 120, 130, 140, 150, 160, 165, 170, 180, 185, 190, 195, 200)
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 25 + $i
    $ws.Range("A$row").Value = "CE" + $values[$i]
}

# Update the sheet view: scroll so row 13 is at the top of the window and
# select A25:A38 (instead of the single cell V26).
$ws.Range("A25:A38").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
